$p = $ppt.ActivePresentation
$s = $p.Slides.Item(23)
$shape = $s.Shapes.Item(2)
$shape.TextFrame.TextRange.Text = "Erstellt Smartmeter und Measurands und nimmt die Testmessungen vor.`r" + "`r" + "SMEmu Restful?"
